$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the designator list for row 7 (swap C2 for C5)
$ws.Range("A7").Value = "C1, C5, C8, C9, C12, C14, C16, C19, C22"

# Correct the manufacturer part number for X4 (row 27)
$ws.Range("B27").Value = "FC-135"

# Correct the manufacturer part number for CHANNEL1, CHANNEL2 (row 5)
$ws.Range("B5").Value = "GS009S-5.0-03P-11"

# Un-bold the designator cells for R17 and C1/C5/... rows to match the rest of column A
$ws.Range("A6").Font.Bold = $false
$ws.Range("A7").Font.Bold = $false

# Move the active selection to B6, matching the saved view state
$ws.Range("B6").Select()
